# This script applies the "Natmi following Dr Hou advice" update to the
# LR-pairs (Fn1-Itga2) results sheet. The underlying NATMI computation was
# re-run, which changed the per-pair statistics for the existing
# Sending/Target cluster combinations (rows 2-7) and added the three
# previously-missing combinations where Sending cluster == Target cluster
# and the "sCs" sending-cluster combinations (new rows 8-10), bringing the
# table to the full 3x3 cross-product of ECs / FAPs / sCs clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-7) with recomputed statistics ---
# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.242496666666667
$ws.Range("N2").Value = 3.72749
$ws.Range("O2").Value = 0.3272238221337332
$ws.Range("P2").Value = 0.3272238221337332
$ws.Range("Q2").Value = 27.14826473577111
$ws.Range("R2").Value = 244.33438262194
$ws.Range("S2").Value = 0.01642669768657148
$ws.Range("T2").Value = 0.01642669768657148

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.744414
$ws.Range("N3").Value = 5.233242000000001
$ws.Range("O3").Value = 0.459408730644692
$ws.Range("P3").Value = 0.459408730644692
$ws.Range("Q3").Value = 38.11504235889467
$ws.Range("R3").Value = 343.0353812300521
$ws.Range("S3").Value = 0.02306240506471345
$ws.Range("T3").Value = 0.02306240506471345

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8101743333333333
$ws.Range("N4").Value = 2.430523
$ws.Range("O4").Value = 0.2133674472215748
$ws.Range("P4").Value = 0.2133674472215748
$ws.Range("Q4").Value = 17.70212176300422
$ws.Range("R4").Value = 159.319095867038
$ws.Range("S4").Value = 0.01071108615751049
$ws.Range("T4").Value = 0.01071108615751049

# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.242496666666667
$ws.Range("N5").Value = 3.72749
$ws.Range("O5").Value = 0.3272238221337332
$ws.Range("P5").Value = 0.3272238221337332
$ws.Range("Q5").Value = 478.4264108809322
$ws.Range("R5").Value = 4305.83769792839
$ws.Range("S5").Value = 0.2894831803543367
$ws.Range("T5").Value = 0.2894831803543367

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.744414
$ws.Range("N6").Value = 5.233242000000001
$ws.Range("O6").Value = 0.459408730644692
$ws.Range("P6").Value = 0.459408730644692
$ws.Range("Q6").Value = 671.6909199840514
$ws.Range("R6").Value = 6045.218279856464
$ws.Range("S6").Value = 0.4064224284233868
$ws.Range("T6").Value = 0.4064224284233868

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8101743333333333
$ws.Range("N7").Value = 2.430523
$ws.Range("O7").Value = 0.2133674472215748
$ws.Range("P7").Value = 0.2133674472215748
$ws.Range("Q7").Value = 311.9596284506615
$ws.Range("R7").Value = 2807.636656055953
$ws.Range("S7").Value = 0.1887585286518176
$ws.Range("T7").Value = 0.1887585286518176

# --- Add new rows (8-10) completing the 3x3 cluster cross-product ---
# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.242496666666667
$ws.Range("N8").Value = 3.72749
$ws.Range("O8").Value = 0.3272238221337332
$ws.Range("P8").Value = 0.3272238221337332
$ws.Range("Q8").Value = 35.22537565590333
$ws.Range("R8").Value = 317.02838090313
$ws.Range("S8").Value = 0.02131394409282503
$ws.Range("T8").Value = 0.02131394409282503

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.744414
$ws.Range("N9").Value = 5.233242000000001
$ws.Range("O9").Value = 0.459408730644692
$ws.Range("P9").Value = 0.459408730644692
$ws.Range("Q9").Value = 49.454972474306
$ws.Range("R9").Value = 445.094752268754
$ws.Range("S9").Value = 0.02992389715659166
$ws.Range("T9").Value = 0.02992389715659166

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8101743333333333
$ws.Range("N10").Value = 2.430523
$ws.Range("O10").Value = 0.2133674472215748
$ws.Range("P10").Value = 0.2133674472215748
$ws.Range("Q10").Value = 22.96883042350566
$ws.Range("R10").Value = 206.719473811551
$ws.Range("S10").Value = 0.01389783241224668
$ws.Range("T10").Value = 0.01389783241224668
